$d = $word.ActiveDocument
$vbreak = [char]11

function Set-ParaText($index, $text) {
    $p = $d.Paragraphs.Item($index)
    $rng = $p.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $text
}

function Set-ParaLines($index, [string[]]$lines) {
    $p = $d.Paragraphs.Item($index)
    $rng = $p.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $lines[0]
    for ($i = 1; $i -lt $lines.Length; $i++) {
        $rng.Collapse(0) | Out-Null
        $rng.InsertAfter($vbreak)
        $rng.Collapse(0) | Out-Null
        $rng.InsertAfter($lines[$i])
    }
}

# --- Word paragraph 6 (0-indexed 5): PT "Objetivos" body -> short "Mercado..." summary line ---
Set-ParaText 6 "Mercado - Estimativa de investimento: - Análise Econômica de Investimentos"

# --- Word paragraph 7 (0-indexed 6): EN "Objetivos" body (italic) -> short "Market..." summary line ---
Set-ParaText 7 "Market - Estimated investment : - Economic Analysis of Investments"

# --- Word paragraph 9 (0-indexed 8): "5840671 - Francisco..." -> PT objectives (3 lines) ---
$ptLine1 = "1) Formativos: Propiciar ao educando as condições básicas e necessárias para a sua formação profissional. "
$ptLine2 = "2) Informativos: fornecer ao educando os conceitos básicos para o entendimento, assessoramento e acompanhamento de Projetos na Indústria Química seguindo metodologia especifica."
$ptLine3 = "3) Automatizantes: desenvolver no educando o raciocínio analítico, obedecendo metodologia sistemática aplicada em projetos."
Set-ParaLines 9 @($ptLine1, $ptLine2, $ptLine3)

# --- Word paragraph 11 (0-indexed 10): short "Mercado - Estimativa..." -> long "Mercado - evolução do mercado..." summary ---
$longPrograma = "Mercado - evolução do mercado - Marketing e análise de mercado - Estimativa de investimento: - capital de giro - capital fixo - bens tangíveis/bens intangíveis - investimentos preliminares/investimentos reais - custos fixos/custos variáveis, depreciação, Conceito econômico de externalidades e abordagens teóricas, Elementos para internalizar as externalidades, Controle direto ou regularização na empresa, métodos indiretos c dados observados, métodos indiretos c dados supostos, métodos diretos c dados supostos, métodos diteros c dados observados, Análise Econômica de Investimentos (aspectos básicos) - Técnicas Estatísticas e taxas de juros - Aspectos básicos de Engenharia Econômica  Distribuição Beta para análise em ambiente de risco."
Set-ParaText 11 $longPrograma

# --- Word paragraph 12 (0-indexed 11): short "Market - Estimated..." (italic) -> EN objectives (3 lines) ---
$enLine1 = "1 ) Formative : Offer the learner the basic conditions necessary for their vocational training."
$enLine2 = "2 ) Informational : provide the student the basics to understanding, advice and monitoring of Projects in the Chemical Industry following specific methodology ."
$enLine3 = "3 ) :The learner develop analytical reasoning , following systematic methodology applied in projects ."
Set-ParaLines 12 @($enLine1, $enLine2, $enLine3)

# --- Word paragraph 17 (0-indexed 16): rotate the three "Avaliação" answers ---
# Use temporary placeholders first to avoid any collision between the rotated values.
$p = $d.Paragraphs.Item(17)
$rng = $p.Range
$rng.Find.Execute("Por meio de aulas presenciais, com apresentação dos fundamentos, e resolução de exercícios e exemplos aplicativos com uso de tabelas e normas específicas.", $false, $true, $false, $false, $false, $true, 1, $false, "__TEMP_METODO_ANSWER__", 2) | Out-Null

$p = $d.Paragraphs.Item(17)
$rng = $p.Range
$rng.Find.Execute("A Avaliação será: (P1 + 2P2)/3 Onde: P1: Prova Individual - c/ peso-1 P2: : Prova Individual - c/ peso-2", $false, $true, $false, $false, $false, $true, 1, $false, "__TEMP_CRITERIO_ANSWER__", 2) | Out-Null

$p = $d.Paragraphs.Item(17)
$rng = $p.Range
$rng.Find.Execute("Prova de exame.", $false, $true, $false, $false, $false, $true, 1, $false, "__TEMP_NORMA_ANSWER__", 2) | Out-Null

$p = $d.Paragraphs.Item(17)
$rng = $p.Range
$rng.Find.Execute("__TEMP_METODO_ANSWER__", $false, $true, $false, $false, $false, $true, 1, $false, "A Avaliação será: (P1 + 2P2)/3 Onde: P1: Prova Individual - c/ peso-1 P2: : Prova Individual - c/ peso-2", 2) | Out-Null

$p = $d.Paragraphs.Item(17)
$rng = $p.Range
$rng.Find.Execute("__TEMP_CRITERIO_ANSWER__", $false, $true, $false, $false, $false, $true, 1, $false, "Prova de exame.", 2) | Out-Null

$p = $d.Paragraphs.Item(17)
$rng = $p.Range
$rng.Find.Execute("__TEMP_NORMA_ANSWER__", $false, $true, $false, $false, $false, $true, 1, $false, "Engenharia econômica e análise de custos. Henrique Hirschfeld. 7 ed. editora atlas. 2007", 2) | Out-Null

# --- Word paragraph 19 (0-indexed 18): Bibliografia entry -> Docente name ---
Set-ParaText 19 "5840671 - Francisco José Moreira Chaves"
